# Sprint12.xlsx - "Added last sprint infos"
#
# Fills in the assignee / daily-burn cells for the tasks that were still
# open on the "Sprint" sheet, and marks every task "Done" now that the
# whole team has reported their numbers for the end of the sprint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# --- Task block owned by Vanja Cvetkovic (rows 5-7): assign her and log
#     the day (column R = "Day 11") each task's remaining effort was burned.
$ws.Range("D5").Value = "Vanja Cvetković"
$ws.Range("R5").Value = 3

$ws.Range("D6").Value = "Vanja Cvetković"
$ws.Range("R6").Value = 5

$ws.Range("D7").Value = "Vanja Cvetković"
$ws.Range("R7").Value = 1

# --- Task block owned by Djuro Nenadovic (rows 8-10): already assigned -
#     log the day each task's effort was actually completed.
$ws.Range("I8").Value = 2
$ws.Range("K9").Value = 3
$ws.Range("K10").Value = 3

# --- Task block owned by Predrag Dimitrijevic (rows 11-13): assign him
#     and log the completion day for each task.
$ws.Range("D11").Value = "Predrag Dimitrijević"
$ws.Range("K11").Value = 5

$ws.Range("D12").Value = "Predrag Dimitrijević"
$ws.Range("M12").Value = 13

$ws.Range("D13").Value = "Predrag Dimitrijević"
$ws.Range("N13").Value = 2

# --- Every task is now finished - set the whole status column to "Done".
$ws.Range("F8").Value = "Done"
$ws.Range("F9").Value = "Done"
$ws.Range("F10").Value = "Done"
$ws.Range("F11").Value = "Done"
$ws.Range("F12").Value = "Done"
$ws.Range("F13").Value = "Done"

# Leave the selection where the author ended up.
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
